$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook theme rename (best effort; engine may not persist names) ---
try {
    $theme = $wb.Theme
    $theme.Name = "Office 2013 - 2022 Theme"
    $cs = $theme.ThemeColorScheme
    $cs.Name = "Office 2013 - 2022"
    $fs = $theme.ThemeFontScheme
    $fs.Name = "Office 2013 - 2022"
} catch { }

# --- Apply AutoFilter on B column (Model) to the 5 "_abundance" values ---
$rng = $ws.Range("A1:I29")
$vals = @("annelid_abundance","crustacea_abundance","ept_abundance","insect_abundance","mollusc_abundance")
$rng.AutoFilter(2, $vals, 7)

# Register the hidden _FilterDatabase defined name at sheet scope (mirrors
# what Excel does automatically when AutoFilter criteria are applied).
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=LT_Yr_metaanaly_weighted_noRand!`$A`$1:`$I`$29")
$fdb.Visible = $false

# The real author's filter left two "_richness" rows (ept_richness / row 16,
# insect_richness / row 20) visible even though they don't match the 5
# selected "_abundance" values -- reproduce that by unhiding them again.
$ws.Rows(16).Hidden = $false
$ws.Rows(20).Hidden = $false

# --- Decrease-decimal formatting on column F (D*(E/100)) for most rows ---
$ws.Range("F2:F15").NumberFormat = "0"
$ws.Range("F17:F19").NumberFormat = "0"
$ws.Range("F21").NumberFormat = "0"
$ws.Range("F23").NumberFormat = "0"
$ws.Range("F25").NumberFormat = "0"
$ws.Range("F27").NumberFormat = "0"

# --- Move the selection / scroll position ---
$ws.Range("C39").Select()

Write-Host "edit applied"
